$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row that only held the "5111420 - Talita Martins Lacerda" value in
# columns B/C (no label in column A) is removed; every row below it shifts
# up by one.
$ws.Rows("13").Delete()

# After the shift, a handful of cells hold text that belongs to a different
# field than before, so refresh those values in place (columns B and C
# always mirror each other in this sheet).
$ws.Range("B10").Value = "5111420 - Talita Martins Lacerda"
$ws.Range("C10").Value = "5111420 - Talita Martins Lacerda"

$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# "01/01/2017" looks like a date, so force text formatting before writing it
# or Excel will silently convert it to a date serial number. Then restore
# the normal body-text style (copied from an existing plain-text cell) so
# no stray date number format sticks to the cell.
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "01/01/2017"
$ws.Range("B3").Copy()
$ws.Range("B15").PasteSpecial(-4122)

$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "01/01/2017"
$ws.Range("C3").Copy()
$ws.Range("C15").PasteSpecial(-4122)

$ws.Range("B18").Value = "5111420 - Talita Martins Lacerda"
$ws.Range("C18").Value = "5111420 - Talita Martins Lacerda"

$ws.Range("B19").Value = "Duas provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."
$ws.Range("C19").Value = "Duas provas escritas envolvendo o conteúdo teórico ministrado em sala de aula."

$ws.Range("B20").Value = "A nota final corresponderá à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados."
$ws.Range("C20").Value = "A nota final corresponderá à média aritmética das duas provas. Os alunos que apresentarem média igual ou superior a 5 estarão aprovados, enquanto que aqueles que tiverem média inferior a 3 estarão reprovados."

$ws.Range("B21").Value = "Será aplicada uma nova avaliação aos alunos com notas finais situadas no intervalo de 3 a 4,9. A nota final do aluno será a média aritmética desta avaliação com a anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
$ws.Range("C21").Value = "Será aplicada uma nova avaliação aos alunos com notas finais situadas no intervalo de 3 a 4,9. A nota final do aluno será a média aritmética desta avaliação com a anteriormente obtida, estando aprovados os alunos que tiverem nota final igual ou superior a 5."
